$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows for a new weekly price report, pushing the
# existing rows 936-1015 down to 938-1017.
$ws.Rows("936:937").Insert()

# Row 936: new "Primera" quality entry
$ws.Range("A936").Value = 8
$ws.Range("B936").Value = "Terminal La Palmera de La Serena"
$ws.Range("C936").Value = "Coquimbo"
$ws.Range("D936").Value = 45223
$ws.Range("E936").Value = 4
$ws.Range("F936").Value = 100112043
$ws.Range("G936").Value = "Pepino ensalada"
$ws.Range("H936").Value = "Sin especificar"
$ws.Range("I936").Value = "Primera"
$ws.Range("J936").Value = 600
$ws.Range("K936").Value = 13000
$ws.Range("L936").Value = 14000
$ws.Range("M936").Value = 13500
$ws.Range("N936").Value = "$/caja 60 unidades"
$ws.Range("O936").Value = "Región de Arica y Parinacota"
$ws.Range("P936").Value = 225
$ws.Range("Q936").Value = 60
$ws.Range("R936").Value = "Hortaliza"

# Row 937: new "Segunda" quality entry
$ws.Range("A937").Value = 8
$ws.Range("B937").Value = "Terminal La Palmera de La Serena"
$ws.Range("C937").Value = "Coquimbo"
$ws.Range("D937").Value = 45223
$ws.Range("E937").Value = 4
$ws.Range("F937").Value = 100112043
$ws.Range("G937").Value = "Pepino ensalada"
$ws.Range("H937").Value = "Sin especificar"
$ws.Range("I937").Value = "Segunda"
$ws.Range("J937").Value = 400
$ws.Range("K937").Value = 9000
$ws.Range("L937").Value = 10000
$ws.Range("M937").Value = 9500
$ws.Range("N937").Value = "$/caja 80 unidades"
$ws.Range("O937").Value = "Región de Arica y Parinacota"
$ws.Range("P937").Value = 119
$ws.Range("Q937").Value = 80
$ws.Range("R937").Value = "Hortaliza"
